$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Investment_Cost")

# --- Column widths (set before other formatting so sizing matches as closely as possible) ---
$ws.Range("B1:F1").ColumnWidth = 14.8
$ws.Range("I1").ColumnWidth = 14.72

# --- Header row (row 1): bold + wrap text on B1:F1, update header text ---
# (text updated before the later new-string cells so the shared-string table
# ends up in the same append order as the target workbook)
$ws.Range("B1:F1").Font.Bold = $true
$ws.Range("B1:F1").WrapText = $true

$ws.Range("B1").Value = "Investment_Cost [Euro/MW or MWh] Value 2020"
$ws.Range("C1").Value = "Investment_Cost [Euro/MW or MWh] Value 2020"
$ws.Range("D1").Value = "Investment_Cost [Euro/MW or MWh] Value 2020"
$ws.Range("E1").Value = "Investment_Cost [Euro/MW or MWh] Value 2020"
$ws.Range("F1").Value = "Investment_Cost [Euro/MW or MWh] Value 2020"

$ws.Rows("1:1").RowHeight = 43.5

# --- Row labels: replace Steam_Plant row with Electrolyzer SOEC, add Electric_Steam_Boiler ---
$ws.Range("A5").Value = "Electrolyzer SOEC"
$ws.Range("A8").Value = "Electric_Steam_Boiler"

# --- Numeric data cells B2:F10, formatted with custom number format ---
$ws.Range("B2:F10").NumberFormat = "#,##0.000000"

$ws.Range("B2").Value = 560000
$ws.Range("C2").Value = 560000
$ws.Range("D2").Value = 380000
$ws.Range("E2").Value = 320000
$ws.Range("F2").Value = 290000

$ws.Range("B3").Value = 1900000
$ws.Range("C3").Value = 1400000
$ws.Range("D3").Value = 875000
$ws.Range("E3").Value = 675000
$ws.Range("F3").Value = 475000

$ws.Range("B4").Value = 1900000
$ws.Range("C4").Value = 1425000
$ws.Range("D4").Value = 950000
$ws.Range("E4").Value = 725000
$ws.Range("F4").Value = 500000

$ws.Range("B5").Value = 2900000
$ws.Range("C5").Value = 2075000
$ws.Range("D5").Value = 1250000
$ws.Range("E5").Value = 1050000
$ws.Range("F5").Value = 800000

$ws.Range("B6").Value = 500000
$ws.Range("C6").Value = 500000
$ws.Range("D6").Value = 500000
$ws.Range("E6").Value = 500000
$ws.Range("F6").Value = 500000

$ws.Range("B7").Value = 1350000
$ws.Range("C7").Value = 1350000
$ws.Range("D7").Value = 1090000
$ws.Range("E7").Value = 960000
$ws.Range("F7").Value = 870000

$ws.Range("B8").Value = 150000
$ws.Range("C8").Value = 145000.00000000003
$ws.Range("D8").Value = 140000
$ws.Range("E8").Value = 135000
$ws.Range("F8").Value = 130000

$ws.Range("B9").Value = 0.00013958682300390843
$ws.Range("C9").Value = 0.00013958682300390843
$ws.Range("D9").Value = 0.00013958682300390843
$ws.Range("E9").Value = 0.00013958682300390843
$ws.Range("F9").Value = 0.00013958682300390843

$ws.Range("B10").Value = 0.121
$ws.Range("C10").Value = 0.17049999999999998
$ws.Range("D10").Value = 0.099
$ws.Range("E10").Value = 0.061
$ws.Range("F10").Value = 0.046

# Extra empty-but-formatted cell
$ws.Range("I7").NumberFormat = "#,##0.000000"

# --- Sheet view / selection ---
$ws.Range("K20").Select() | Out-Null

# --- Page setup (orientation) ---
$ws.PageSetup.Orientation = 1

Write-Host "done"
